$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$edits = @(
    ,@(2, 4, '329.08')
    ,@(2, 5, '0.48%')
    ,@(2, 7, '14')
    ,@(3, 5, '0.95%')
    ,@(3, 7, '14')
    ,@(4, 5, '-0.77%')
    ,@(4, 7, '14')
    ,@(5, 4, '0.08009')
    ,@(5, 5, '-0.44%')
    ,@(5, 7, '14')
    ,@(6, 4, '2.066')
    ,@(6, 5, '8.65%')
    ,@(6, 7, '14')
    ,@(7, 2, 'MXToken')
    ,@(7, 3, 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx')
    ,@(7, 4, '0.9547')
    ,@(7, 5, '1.06%')
    ,@(7, 7, '14')
    ,@(8, 2, 'LiechtensteinCryptoassetsExchange')
    ,@(8, 3, 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx')
    ,@(8, 4, '0.1140')
    ,@(8, 5, '-1.22%')
    ,@(8, 7, '14')
    ,@(9, 2, 'WazirX')
    ,@(9, 3, 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx')
    ,@(9, 4, '0.1886')
    ,@(9, 5, '2.64%')
    ,@(9, 7, '14')
    ,@(10, 2, 'MCDex')
    ,@(10, 3, 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb')
    ,@(10, 4, '10.17')
    ,@(10, 5, '5.96%')
    ,@(10, 7, '14')
    ,@(11, 2, 'MandalaExchangeToken')
    ,@(11, 3, 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx')
    ,@(11, 4, '0.09892')
    ,@(11, 5, '1.79%')
    ,@(11, 7, '14')
    ,@(12, 2, 'BitrueCoin')
    ,@(12, 3, 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr')
    ,@(12, 4, '0.04873')
    ,@(12, 5, '11.14%')
    ,@(12, 7, '14')
    ,@(13, 2, 'BitMartToken')
    ,@(13, 3, 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx')
    ,@(13, 4, '0.1062')
    ,@(13, 5, '-0.50%')
    ,@(13, 7, '14')
    ,@(14, 2, 'BitForexToken')
    ,@(14, 3, 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf')
    ,@(14, 4, '0.001270')
    ,@(14, 5, '-0.90%')
    ,@(14, 7, '14')
    ,@(15, 2, 'CoinExToken')
    ,@(15, 3, 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet')
    ,@(15, 4, '0.04082')
    ,@(15, 5, '-3.12%')
    ,@(15, 7, '14')
    ,@(16, 2, 'TigerCash')
    ,@(16, 3, 'https://coinranking.com/coin/6hIn06L2+tigercash-tch')
    ,@(16, 4, '0.006148')
    ,@(16, 5, '2.87%')
    ,@(16, 7, '14')
    ,@(17, 2, 'LEO')
    ,@(17, 3, 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo')
    ,@(17, 4, '3.378')
    ,@(17, 5, '-0.84%')
    ,@(17, 7, '14')
    ,@(18, 2, 'GateToken')
    ,@(18, 3, 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt')
    ,@(18, 4, '4.413')
    ,@(18, 5, '3.09%')
    ,@(18, 7, '14')
    ,@(19, 2, 'BTSEToken')
    ,@(19, 3, 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse')
    ,@(19, 4, '2.627')
    ,@(19, 5, '2.35%')
    ,@(19, 7, '14')
    ,@(20, 4, '0.3407')
    ,@(20, 5, '-2.42%')
    ,@(20, 7, '14')
    ,@(21, 5, '-0.53%')
    ,@(21, 7, '14')
    ,@(22, 5, '2.78%')
    ,@(22, 7, '14')
    ,@(23, 4, '0.001301')
    ,@(23, 5, '4.19%')
    ,@(23, 7, '14')
    ,@(24, 4, '0.004351')
    ,@(24, 5, '1.18%')
    ,@(24, 7, '14')
    ,@(25, 4, '0.0001251')
    ,@(25, 5, '-0.90%')
    ,@(25, 7, '14')
    ,@(26, 4, '0.0003746')
    ,@(26, 5, '-6.28%')
    ,@(26, 7, '14')
    ,@(27, 7, '14')
    ,@(28, 7, '14')
    ,@(29, 7, '14')
    ,@(30, 7, '14')
    ,@(31, 7, '14')
    ,@(32, 7, '14')
    ,@(33, 7, '14')
    ,@(34, 7, '14')
    ,@(35, 7, '14')
    ,@(36, 7, '14')
    ,@(37, 7, '14')
    ,@(38, 4, '0.02581')
    ,@(38, 5, '-2.31%')
    ,@(38, 7, '14')
    ,@(39, 4, '0.05773')
    ,@(39, 5, '5.29%')
    ,@(39, 7, '14')
    ,@(40, 4, '0.007593')
    ,@(40, 5, '0.17%')
    ,@(40, 7, '14')
    ,@(41, 4, '0.1401')
    ,@(41, 5, '0.39%')
    ,@(41, 7, '14')
    ,@(42, 4, '0.007320')
    ,@(42, 5, '-9.63%')
    ,@(42, 7, '14')
    ,@(43, 4, '0.002009')
    ,@(43, 5, '-0.14%')
    ,@(43, 7, '14')
    ,@(44, 4, '0.009074')
    ,@(44, 5, '2.53%')
    ,@(44, 7, '14')
    ,@(45, 5, '1.20%')
    ,@(45, 7, '14')
    ,@(46, 5, '-0.15%')
    ,@(46, 7, '14')
    ,@(47, 4, '0.0005803')
    ,@(47, 5, '-0.14%')
    ,@(47, 7, '14')
    ,@(48, 4, '0.003500')
    ,@(48, 5, '53.87%')
    ,@(48, 7, '14')
    ,@(49, 5, '-1.13%')
    ,@(49, 7, '14')
    ,@(50, 5, '-0.15%')
    ,@(50, 7, '14')
    ,@(51, 5, '-0.15%')
    ,@(51, 7, '14')
)

foreach ($edit in $edits) {
    $row = $edit[0]
    $col = $edit[1]
    $val = $edit[2]
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}
